$d = $word.ActiveDocument

# Move to the very end of the document content and append a line break
# (an empty <w:r><w:br/></w:r> run), matching a Shift+Enter press with
# nothing typed afterwards yet.
$rng = $d.Range($d.Content.End, $d.Content.End)
$rng.InsertBreak(6)   # wdLineBreak

# Append a second line break followed immediately by the new text, as if
# the user pressed Shift+Enter again and then typed the sentence (both
# land together in one run, just like the previous "Segunda alteração do
# 002" line).
$rng = $d.Range($d.Content.End, $d.Content.End)
$rng.InsertAfter([char]11 + "2** Acho que essa modificação é correta")
$endPos = $rng.End

# Word leaves a "_GoBack" bookmark at the last edited location after
# typing. A zero-length bookmark placed exactly at the end of a
# paragraph's text (immediately before its paragraph mark) normalizes to
# span the whole paragraph, so temporarily add a placeholder character
# after the insertion point, anchor the bookmark there, and then remove
# the placeholder again -- leaving a proper collapsed bookmark right
# after the last run.
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Delete()
